$wb = $excel.ActiveWorkbook

# The long diagnostic message now shown in the "Error Detail" column (P) for
# the handback row, for both the zh-cn and de-de sheets.
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d06bd74ac638215f09aa67bc15bce132c5cb83d3/e2e/7958228f-0342-4475-b196-c8e10920b047.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4707a95597912fe09953308c8049f162bd1abfe1/e2e/7958228f-0342-4475-b196-c8e10920b047.md."

# The "latest" target-file hyperlink address, reused for the new I5 hyperlinks
# on both locale sheets (matches the existing A5 "current" hyperlink target).
$latestHandoffUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4707a95597912fe09953308c8049f162bd1abfe1/e2e/7958228f-0342-4475-b196-c8e10920b047.md"

function Update-LocaleSheet {
    param(
        [string]$SheetName,
        [string]$HandbackFile,
        [string]$HandbackDateTime
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Widen the "Error Detail" column (P / column 16) so the long message is
    # readable.
    $ws.Columns.Item(16).ColumnWidth = 39.17

    # I5 - "Latest Target File": now links to the handback markdown file.
    $ws.Hyperlinks.Add($ws.Range("I5"), $latestHandoffUrl, "", "", "7958228f-0342-4475-b196-c8e10920b047.md") | Out-Null
    $ws.Range("I5").Font.Underline = $true
    $ws.Range("I5").Font.Color = 15570276

    # J5 - "Latest Handback File": the generated xliff file name for this locale.
    $ws.Range("J5").Value2 = $HandbackFile

    # K5 - "Latest Handback DateTime": timestamp of the handback generation.
    $ws.Range("K5").Value2 = $HandbackDateTime

    # P5 - "Error Detail": the handback-is-stale diagnostic.
    $ws.Range("P5").Value2 = $errorDetail
}

Update-LocaleSheet -SheetName "zh-cn" -HandbackFile "7958228f-0342-4475-b196-c8e10920b047.323c40e87ff7c4d66c22f4cf1f91bf2ae148b175.zh-cn.xlf" -HandbackDateTime "2016-09-09 12:33:00"
Update-LocaleSheet -SheetName "de-de" -HandbackFile "7958228f-0342-4475-b196-c8e10920b047.323c40e87ff7c4d66c22f4cf1f91bf2ae148b175.de-de.xlf" -HandbackDateTime "2016-09-09 12:33:25"
